$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.390.66'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '2.455.72'
$ws.Range('E3').Value = '  +8.06%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '297.47'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.00%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '97.43'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -2.60%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.578'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +0.75%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +0.06%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.516'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +1.35%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '35.56'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +1.03%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0789'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -1.63%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '7.20'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +1.93%  '
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('D14').Value = '2.827.40'
$ws.Range('E14').Value = '  +7.96%  '
$ws.Range('D15').Value = '2.446.01'
$ws.Range('E15').Value = '  +7.51%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.859'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +7.36%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '14.16'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +3.76%  '
$ws.Range('D18').Value = '46.385.41'
$ws.Range('E18').Value = '  -0.43%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '12.92'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +3.01%  '
$ws.Range('D20').Value = '0.0₃0953'
$ws.Range('E20').Value = '  -3.96%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.29'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +7.40%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '67.80'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +2.77%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '247.00'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('E24').Value = '  +0.75%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '1.98'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +6.10%  '
$ws.Range('E26').Value = '  -0.08%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '40.04'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -3.21%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -1.29%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '9.87'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +2.67%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '3.87'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +15.04%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '21.56'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +6.95%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '5.64'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +5.21%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -1.58%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '148.77'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +1.30%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.0778'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +0.88%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +21.86%  '
$ws.Range('E37').Value = '  +1.51%  '
$ws.Range('E38').Value = '  +0.22%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '15.47'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -1.43%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '3.96'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +2.32%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.0304'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +2.14%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '3.32'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +7.23%  '
$ws.Range('D43').Value = '1.992.46'
$ws.Range('E43').Value = '  +11.56%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +0.05%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '92.79'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('E46').Value = '  -2.60%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '16.61'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +33.54%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '8.61'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +8.86%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '102.29'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +7.53%  '
$ws.Range('D50').Value = '2.693.65'
$ws.Range('E50').Value = '  +7.98%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.188'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +1.44%  '
